# repull data, push all data, mean calculation
# Update column F (dSF) values for the specific rows that were repulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 4
    9  = 3
    14 = 0
    16 = 1
    18 = -1
    19 = 2
    21 = 1
    33 = -5
    34 = -3
    35 = 1
    36 = 0
    37 = -2
    40 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
